# Convert the four Word field codes ({m:comment ...}, {m:userdoc 'zone1'},
# {m:enduserdoc}) that are currently stored as real Word fields
# (fldChar begin/instrText/fldChar end) into plain literal-text runs, e.g.
# "{m:comment Title1}", while keeping the run that carries the heading
# character style ("Title1" / "Title2") as its own separate run.
#
# Strategy per field:
#   1. Remember the paragraph that hosts the field and the field's code
#      text (used to decide exactly what literal text to emit).
#   2. Delete the field (removes the begin/instrText*/end runs cleanly).
#   3. Use Range.InsertXML to drop in freshly built <w:r><w:t>...</w:t></w:r>
#      runs with the exact text (and xml:space="preserve" where needed) -
#      InsertXML lets us control whitespace-preservation explicitly,
#      unlike plain Range.Text assignment which infers it from content.
#   4. Re-apply the character style (e.g. Titre1Car) to the title
#      sub-string via Range.Style, since InsertXML does not carry rStyle.

$d = $word.ActiveDocument

$wOpenXmlNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function New-RunsPackageXml($runsXml) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document ' + $wOpenXmlNS + '><w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

function Set-LiteralRun($t, $preserve) {
    if ($preserve) {
        return '<w:r><w:t xml:space="preserve">' + $t + '</w:t></w:r>'
    } else {
        return '<w:r><w:t>' + $t + '</w:t></w:r>'
    }
}

# --- Paragraph 1: {m:comment Title1} -------------------------------------
$p = $d.Paragraphs.Item(1)
$d.Fields.Item(1).Delete()
$r = $p.Range
$pStart = $r.Start
$runs = (Set-LiteralRun "{m:comment " $true) + (Set-LiteralRun "Title1" $false) + (Set-LiteralRun "}" $true)
$r.InsertXML((New-RunsPackageXml $runs))
$titleStart = $pStart + 11
$titleEnd = $titleStart + 6
$d.Range($titleStart, $titleEnd).Style = "Titre1Car"

# --- Paragraph 2: {m:comment Title2} -------------------------------------
$p = $d.Paragraphs.Item(2)
$d.Fields.Item(1).Delete()
$r = $p.Range
$pStart = $r.Start
$runs = (Set-LiteralRun "{m:comment " $true) + (Set-LiteralRun "Title2" $false) + (Set-LiteralRun "}" $true)
$r.InsertXML((New-RunsPackageXml $runs))
$titleStart = $pStart + 11
$titleEnd = $titleStart + 6
$d.Range($titleStart, $titleEnd).Style = "Titre2Car"

# --- Paragraph 4: {m:userdoc 'zone1'} -------------------------------------
$p = $d.Paragraphs.Item(4)
$d.Fields.Item(1).Delete()
$r = $p.Range
$runs = (Set-LiteralRun "{" $false) + (Set-LiteralRun "m" $false) + (Set-LiteralRun ":userdoc 'zone1'" $false) + (Set-LiteralRun "}" $true)
$r.InsertXML((New-RunsPackageXml $runs))

# --- Paragraph 6: {m:enduserdoc} ------------------------------------------
$p = $d.Paragraphs.Item(6)
$d.Fields.Item(1).Delete()
$r = $p.Range
$runs = Set-LiteralRun "{m:enduserdoc}" $true
$r.InsertXML((New-RunsPackageXml $runs))

Write-Host "Done"
